# Scheduled-runner update: refresh Universalis market-price-derived
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) across the
# Zalera_Profits crafting-leve sheets. Values mirror the latest market
# snapshot; a handful of LeveProfitNQ/LeveProfitHQ cells that no longer
# have a valid comparison (HQ leve price without an NQ counterpart, etc.)
# are cleared rather than recomputed.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 536.62
$ws.Range("J17").Value = 536.62
$ws.Range("L17").Value = 1609.86
$ws.Range("N17").Value = -1945.86
# row 18
$ws.Range("H18").Value = 1111.7273
$ws.Range("I18").Value = 972.9
$ws.Range("K18").Value = 972.9
$ws.Range("M18").Value = -688.9
# row 40
$ws.Range("H40").Value = 2599.0557
$ws.Range("I40").Value = 3436.1428
$ws.Range("J40").Value = 2066.3635
$ws.Range("K40").Value = 3436.1428
$ws.Range("L40").Value = 2066.3635
$ws.Range("M40").Value = -3261.1428
$ws.Range("N40").Value = -2416.3635
# row 62
$ws.Range("H62").Value = 79167940
$ws.Range("I62").Value = 90477070
$ws.Range("J62").Value = 3999.5
$ws.Range("K62").Value = 90477070
$ws.Range("L62").Value = 3999.5
$ws.Range("M62").Value = -90476446
$ws.Range("N62").Value = -5247.5
# row 65
$ws.Range("H65").Value = 79167940
$ws.Range("I65").Value = 90477070
$ws.Range("J65").Value = 3999.5
$ws.Range("K65").Value = 452385350
$ws.Range("L65").Value = 19997.5
$ws.Range("M65").Value = -452382230
$ws.Range("N65").Value = -26237.5
# row 98
$ws.Range("H98").Value = 1052.2727
$ws.Range("I98").Value = 1062.5
$ws.Range("K98").Value = 1062.5
$ws.Range("M98").Value = 435.5
# row 100
$ws.Range("H100").Value = 2664.6365
$ws.Range("J100").Value = 2761.1
$ws.Range("L100").Value = 2761.1
$ws.Range("N100").Value = -3843.1
# row 111
$ws.Range("H111").Value = 2515.4443
$ws.Range("I111").Value = 2451
$ws.Range("K111").Value = 7353
$ws.Range("M111").Value = -4286
# row 122
$ws.Range("H122").Value = 1052.2727
$ws.Range("I122").Value = 1062.5
$ws.Range("K122").Value = 3187.5
$ws.Range("M122").Value = -737.5
# row 137
$ws.Range("H137").Value = 6586968
$ws.Range("I137").Value = 11908646
$ws.Range("J137").Value = 13129.647
$ws.Range("K137").Value = 35725938
$ws.Range("L137").Value = 39388.94100000001
$ws.Range("M137").Value = -35723388
$ws.Range("N137").Value = -44488.94100000001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 132
$ws.Range("H132").Value = 6166.0835
$ws.Range("I132").Value = 4165.8887
$ws.Range("K132").Value = 12497.6661
$ws.Range("M132").Value = -9967.666100000002
# row 134
$ws.Range("H134").Value = 83999
$ws.Range("J134").Value = 83999
$ws.Range("L134").Value = 83999
$ws.Range("N134").Value = -94139

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 57
$ws.Range("H57").Value = 97999.336
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 97999.336
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 97999.336
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -99439.336
# row 94
$ws.Range("H94").Value = 1183.8
$ws.Range("I94").Value = 1229.8334
$ws.Range("K94").Value = 1229.8334
$ws.Range("M94").Value = -778.8334
# row 132
$ws.Range("H132").Value = 99991.336
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# row 133
$ws.Range("H133").Value = 119999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 119999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 119999
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -130119
# row 134
$ws.Range("H134").Value = 7408.3125
$ws.Range("I134").Value = 2040
$ws.Range("K134").Value = 6120
$ws.Range("M134").Value = -3585
# row 136
$ws.Range("H136").Value = 97999.336
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 97999.336
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 97999.336
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -108199.336

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 28
$ws.Range("H28").Value = 34750
$ws.Range("J28").Value = 34750
$ws.Range("L28").Value = 34750
$ws.Range("N28").Value = -35240
# row 122
$ws.Range("H122").Value = 126112
$ws.Range("I122").Value = 200799.4
$ws.Range("J122").Value = 1633
$ws.Range("K122").Value = 602398.2
$ws.Range("L122").Value = 4899
$ws.Range("M122").Value = -599948.2
$ws.Range("N122").Value = -9799

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 3801949.2
$ws.Range("I4").Value = 2735100.5
$ws.Range("J4").Value = 5443255.5
$ws.Range("K4").Value = 8205301.5
$ws.Range("L4").Value = 16329766.5
$ws.Range("M4").Value = -8205189.5
$ws.Range("N4").Value = -16329990.5
# row 33
$ws.Range("H33").Value = 53.708332
$ws.Range("I33").Value = 53.708332
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 322.249992
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -39.24999200000002
$ws.Range("N33").ClearContents()
# row 92
$ws.Range("H92").Value = 2752.6843
$ws.Range("J92").Value = 4300.1665
$ws.Range("L92").Value = 12900.4995
$ws.Range("N92").Value = -15396.4995
# row 140
$ws.Range("H140").Value = 53772.844
$ws.Range("I140").Value = 77898.766
$ws.Range("K140").Value = 233696.298
$ws.Range("M140").Value = -228516.298

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 1325
$ws.Range("I2").Value = 1868.8334
$ws.Range("J2").Value = 237.33333
$ws.Range("K2").Value = 1868.8334
$ws.Range("L2").Value = 237.33333
$ws.Range("M2").Value = -1755.8334
$ws.Range("N2").Value = -463.33333
# row 132
$ws.Range("H132").Value = 8270.888999999999
$ws.Range("I132").Value = 4287.6
$ws.Range("J132").Value = 13250
$ws.Range("K132").Value = 12862.8
$ws.Range("L132").Value = 39750
$ws.Range("M132").Value = -10332.8
$ws.Range("N132").Value = -44810

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 20500
# row 71
$ws.Range("H71").Value = 20500
# row 93
$ws.Range("H93").Value = 896.8182
$ws.Range("J93").Value = 1000
$ws.Range("L93").Value = 1000
$ws.Range("N93").Value = -3496
# row 132
$ws.Range("H132").Value = 8496.5
$ws.Range("I132").Value = 6893.4
$ws.Range("J132").Value = 10099.6
$ws.Range("K132").Value = 20680.2
$ws.Range("L132").Value = 30298.8
$ws.Range("M132").Value = -18150.2
$ws.Range("N132").Value = -35358.8

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 4095.2354
$ws.Range("I132").Value = 1918.2222
$ws.Range("J132").Value = 6544.375
$ws.Range("K132").Value = 5754.6666
$ws.Range("L132").Value = 19633.125
$ws.Range("M132").Value = -3224.6666
$ws.Range("N132").Value = -24693.125
